$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.698.84"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.89"
$ws.Range("E3").Value = "  +1.56%  "

$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.34"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4634"
$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3913"
$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07931"
$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9703"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.29"
$ws.Range("E11").Value = "  +1.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.886.58"
$ws.Range("E12").Value = "  +1.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.737"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.943"
$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06978"
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.29"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001008"
$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.97"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.716.35"
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.331"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  +0.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.127"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.154.93"
$ws.Range("E25").Value = "  +4.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.64"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.37"
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.705"
$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.002"
$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.55"
$ws.Range("E30").Value = "  +2.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09377"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9316"
$ws.Range("E32").Value = "  -1.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.332"
$ws.Range("E33").Value = "  +0.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.348"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.360"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05845"
$ws.Range("E36").Value = "  -2.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02131"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.152"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.911"
$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5665"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.944"
$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1785"
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07241"
$ws.Range("E43").Value = "  +2.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.72"
$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5327"
$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.152"
$ws.Range("E46").Value = "  -5.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.140"
$ws.Range("E47").Value = "  -6.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.848"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.65"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.006"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.345"
$ws.Range("E51").Value = "  +0.95%  "
